# Auto-update draw results: append the 2025-11-01 Pick 4 draw as a new
# row (46) at the bottom of the "Results" sheet, and extend the sheet's
# dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 46

# Columns A and C contain values that look like dates / pure numbers
# ("2025-11-01" and "251101"). Excel's automatic type inference would
# otherwise silently convert them into a date serial number / numeric
# value. Prefixing with a leading apostrophe forces them to be stored
# as literal text (matching the existing rows, which are all text).
# Re-applying the "Normal" style afterwards clears the "quote prefix"
# formatting flag that gets attached as a side effect, so the new
# cells stay unstyled just like the rest of the sheet.

$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Value = "'2025-11-01"
$cellA.Style = "Normal"

$cellB = $ws.Cells.Item($newRow, 2)
$cellB.Value = "Pick 4"

$cellC = $ws.Cells.Item($newRow, 3)
$cellC.Value = "'251101"
$cellC.Style = "Normal"

$cellD = $ws.Cells.Item($newRow, 4)
$cellD.Value = "6-5-6-5"

$cellE = $ws.Cells.Item($newRow, 5)
$cellE.Value = "2025-11-01T21:35:28.396+04:00"
